# Add 2022-Q3 data: a new summary row on "总计" and a new detail worksheet
# "2022-Q3" inserted right after it (all later sheets shift right by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" (sheet 1) — insert a new row 2 for the 2022-Q3 summary
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a blank row above the current row 2 (old row 2 becomes row 3, etc).
$ws1.Rows.Item(2).Insert()

# The insert leaves row 2 with a stray copied format — reset it, then copy
# the (untouched) "A column" style from row 3 so A2 matches A3..A7.
$ws1.Range("A2:D2").ClearFormats()
$ws1.Range("A3").Copy()
$ws1.Range("A2").PasteSpecial(-4122)

$ws1.Cells.Item(2, 1).Value = 0
$ws1.Cells.Item(2, 2).Value = "2022-Q3"
$ws1.Cells.Item(2, 3).Value = 3
$ws1.Cells.Item(2, 4).Value = 0.1

# The index column (A) for every pre-existing row shifts up by one.
$ws1.Cells.Item(3, 1).Value = 1
$ws1.Cells.Item(4, 1).Value = 2
$ws1.Cells.Item(5, 1).Value = 3
$ws1.Cells.Item(6, 1).Value = 4
$ws1.Cells.Item(7, 1).Value = 5

# ---------------------------------------------------------------------
# Step 2: insert the new "2022-Q3" worksheet right after "总计" — it
# becomes the new sheet #2, pushing every other quarter sheet back by one.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "2022-Q3"
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# "2022-Q1" (the old sheet 2) is now sheet 3 — same column layout, reuse it
# as a formatting template so fonts/borders/alignment match exactly.
$template = $wb.Worksheets.Item(3)

# ---- header row ----
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Column A (index) style matches the template's A2:A4.
$template.Range("A2:A4").Copy()
$newSheet.Range("A2:A4").PasteSpecial(-4122)

# Data columns B..G hold numeric-looking strings ("4.19", "501021", ...)
# that must round-trip as text (keeps the leading zero in "006127"), so
# force a text format *before* writing the values.
$newSheet.Range("B2:G4").NumberFormat = "@"

# ---- data rows: column A / H are numeric, B..G are text ----
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "501021"
$newSheet.Cells.Item(2, 3).Value = "华宝标普香港上市中国中小盘指数（LOF）A"
$newSheet.Cells.Item(2, 4).Value = "4.19"
$newSheet.Cells.Item(2, 5).Value = "92.99"
$newSheet.Cells.Item(2, 6).Value = "1.72"
$newSheet.Cells.Item(2, 7).Value = "0.0721"
$newSheet.Cells.Item(2, 8).Value = 8

$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "513160"
$newSheet.Cells.Item(3, 3).Value = "银华恒生港股通中国科技ETF"
$newSheet.Cells.Item(3, 4).Value = "0.43"
$newSheet.Cells.Item(3, 5).Value = "92.45"
$newSheet.Cells.Item(3, 6).Value = "5.75"
$newSheet.Cells.Item(3, 7).Value = "0.0247"
$newSheet.Cells.Item(3, 8).Value = 7

$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(4, 2).Value = "006127"
$newSheet.Cells.Item(4, 3).Value = "华宝标普香港上市中国中小盘指数（LOF）C"
$newSheet.Cells.Item(4, 4).Value = "0.24"
$newSheet.Cells.Item(4, 5).Value = "92.99"
$newSheet.Cells.Item(4, 6).Value = "1.72"
$newSheet.Cells.Item(4, 7).Value = "0.0041"
$newSheet.Cells.Item(4, 8).Value = 8

Write-Output "2022-Q3 sheet + summary row added"
